$d = $word.ActiveDocument

# Change 1: date "9/18/2020" -> "8/18/2020" split into two runs "8" and "/18/2020"
$d.Content.Find.Execute("9/18/2020", $false, $false, $false, $false, $false, $true, 1, $false, "8/18/2020", 2)

